# Update master to output generated at c8c62b6
#
# Set each value directly via Cell.Range.Text rather than Find/Replace:
# several of the division problems are textual duplicates of each other
# (e.g. "57÷4=14, 1" appears twice) but map to different replacement
# values depending on which table cell they're in, so a document-wide
# Find.Execute(..., Replace:=wdReplaceAll) would be ambiguous/incorrect.
# Scoping a Range to a single Table.Cell() and assigning .Text directly
# targets only that cell, regardless of duplicate text elsewhere.

$d = $word.ActiveDocument

# Update the date paragraph
$d.Paragraphs.Item(1).Range.Text = "2025-10-22 Wednesday"

$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "78÷5=15, 3"
$t.Cell(1,2).Range.Text = "32÷8=4, 0"
$t.Cell(1,3).Range.Text = "68÷8=8, 4"
$t.Cell(1,4).Range.Text = "91÷2=45, 1"
$t.Cell(1,5).Range.Text = "80÷5=16, 0"

$t.Cell(5,1).Range.Text = "50÷5=10, 0"
$t.Cell(5,2).Range.Text = "80÷4=20, 0"
$t.Cell(5,3).Range.Text = "29÷9=3, 2"
$t.Cell(5,4).Range.Text = "72÷4=18, 0"
$t.Cell(5,5).Range.Text = "13÷7=1, 6"

$t.Cell(9,1).Range.Text = "49÷6=8, 1"
$t.Cell(9,2).Range.Text = "18÷3=6, 0"
$t.Cell(9,3).Range.Text = "12÷3=4, 0"
$t.Cell(9,4).Range.Text = "95÷9=10, 5"
$t.Cell(9,5).Range.Text = "54÷5=10, 4"

$t.Cell(13,1).Range.Text = "97÷2=48, 1"
$t.Cell(13,2).Range.Text = "36÷5=7, 1"
$t.Cell(13,3).Range.Text = "96÷2=48, 0"
$t.Cell(13,4).Range.Text = "17÷3=5, 2"
$t.Cell(13,5).Range.Text = "24÷9=2, 6"

$t.Cell(17,1).Range.Text = "22÷6=3, 4"
$t.Cell(17,2).Range.Text = "25÷3=8, 1"
$t.Cell(17,3).Range.Text = "61÷2=30, 1"
$t.Cell(17,4).Range.Text = "89÷5=17, 4"
$t.Cell(17,5).Range.Text = "67÷4=16, 3"
